# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a completed
# handback: the status text changes from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language sheets (zh-cn, de-de)
# get their "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (with a hyperlink on the target-file cell), and
# a few columns are widened to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Known external targets (same commit as the existing hyperlinks in column A)
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f558e9ae473ef0782bbf1753a6be0ab1bebe4915/e2e/79ffca11-d218-4785-9235-c29d0339b29f.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f558e9ae473ef0782bbf1753a6be0ab1bebe4915/e2e/7b119974-65d2-4413-9d28-4b6615c1c04d.md"

$mdName1 = "79ffca11-d218-4785-9235-c29d0339b29f.md"
$mdName2 = "7b119974-65d2-4413-9d28-4b6615c1c04d.md"

# ---------------------------------------------------------------------------
# 1. Overview sheet: update Status-like columns (E/F) and widen them
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: Status column, Target/Handback file + datetime columns
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Row 2 (79ffca11...)
$zhcn.Range("I2").Value = $mdName1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1) | Out-Null
$zhcn.Range("J2").Value = "79ffca11-d218-4785-9235-c29d0339b29f.55d704c33d5b1872bd722e72ecca78f735b2bb2a.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-25 16:50:29"

# Row 3 (7b119974...)
$zhcn.Range("I3").Value = $mdName2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2) | Out-Null
$zhcn.Range("J3").Value = "7b119974-65d2-4413-9d28-4b6615c1c04d.bf461714b7b9a4d36eadc07c883a904cf7c52506.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-25 16:50:29"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. de-de sheet: Status column, Target/Handback file + datetime columns
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Row 2 (79ffca11...)
$dede.Range("I2").Value = $mdName1
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1) | Out-Null
$dede.Range("J2").Value = "79ffca11-d218-4785-9235-c29d0339b29f.55d704c33d5b1872bd722e72ecca78f735b2bb2a.de-de.xlf"
$dede.Range("K2").Value = "2016-08-25 16:50:36"

# Row 3 (7b119974...)
$dede.Range("I3").Value = $mdName2
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2) | Out-Null
$dede.Range("J3").Value = "7b119974-65d2-4413-9d28-4b6615c1c04d.bf461714b7b9a4d36eadc07c883a904cf7c52506.de-de.xlf"
$dede.Range("K3").Value = "2016-08-25 16:50:36"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

Write-Output "Handback report generated."
